# The source data for "Hortaliza, Vega Modelo de Temuco - Acelga" gained a
# new weekly record. It is inserted as a new row 410 (pushing the existing
# rows 410-482 down to 411-483), matching the diff's row-shift pattern and
# the new dimension A1:R483.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 410; everything below shifts down by one.
$ws.Rows(410).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Cells.Item(410, 1).Value = 10
$ws.Cells.Item(410, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(410, 3).Value = "La Araucanía"
$ws.Cells.Item(410, 4).Value = 45015
$ws.Cells.Item(410, 5).Value = 9
$ws.Cells.Item(410, 6).Value = 100112009
$ws.Cells.Item(410, 7).Value = "Acelga"
$ws.Cells.Item(410, 8).Value = "Sin especificar"
$ws.Cells.Item(410, 9).Value = "Primera"
$ws.Cells.Item(410, 10).Value = 50
$ws.Cells.Item(410, 11).Value = 8000
$ws.Cells.Item(410, 12).Value = 8000
$ws.Cells.Item(410, 13).Value = 8000
$ws.Cells.Item(410, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(410, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(410, 16).Value = 667
$ws.Cells.Item(410, 17).Value = 12
$ws.Cells.Item(410, 18).Value = "Hortaliza"
